# Swap columns A and B (category name <-> item name) on the active sheet,
# for the header row and all data rows, per the commit:
# "修改了sumsales_year下各个分类的列的前后顺序"
# (swapped the column order of category/item name columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aVal = $aCell.Value2
    $bVal = $bCell.Value2

    $aCell.Value = $bVal
    $bCell.Value = $aVal
}
